$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = '@'
    $c.Value = $text
    $c.Style = 'Normal'
}

Set-TextCell 'D2' '61.109.08'
Set-TextCell 'E2' '  +0.36%  '
Set-TextCell 'D3' '2.924.12'
Set-TextCell 'E3' '  +0.18%  '
Set-TextCell 'E4' '  +0.00%  '
Set-TextCell 'D5' '595.78'
Set-TextCell 'E5' '  +0.36%  '
Set-TextCell 'D6' '143.70'
Set-TextCell 'E6' '  -1.28%  '
Set-TextCell 'E7' '  -0.01%  '
Set-TextCell 'D8' '0.500'
Set-TextCell 'E8' '  -1.29%  '
Set-TextCell 'D9' '6.94'
Set-TextCell 'E9' '  +2.14%  '
Set-TextCell 'E10' '  -1.86%  '
Set-TextCell 'D11' '0.437'
Set-TextCell 'E11' '  -0.49%  '
Set-TextCell 'D12' '0.0000223'
Set-TextCell 'E12' '  -1.00%  '
Set-TextCell 'D13' '33.19'
Set-TextCell 'E13' '  -1.23%  '
Set-TextCell 'E14' '  +0.09%  '
Set-TextCell 'D15' '3.406.71'
Set-TextCell 'E15' '  +0.16%  '
Set-TextCell 'D16' '61.133.54'
Set-TextCell 'E16' '  +0.41%  '
Set-TextCell 'D17' '2.920.51'
Set-TextCell 'E17' '  +0.37%  '
Set-TextCell 'D18' '6.64'
Set-TextCell 'E18' '  -0.81%  '
Set-TextCell 'D19' '432.72'
Set-TextCell 'D20' '13.50'
Set-TextCell 'E20' '  +1.08%  '
Set-TextCell 'D21' '0.672'
Set-TextCell 'E21' '  -1.21%  '
Set-TextCell 'D22' '7.04'
Set-TextCell 'E22' '  -0.38%  '
Set-TextCell 'D23' '81.60'
Set-TextCell 'E23' '  +0.05%  '
Set-TextCell 'D24' '10.83'
Set-TextCell 'E24' '  -1.03%  '
Set-TextCell 'E25' '  -2.02%  '
Set-TextCell 'D26' '11.71'
Set-TextCell 'E26' '  -1.73%  '
Set-TextCell 'E27' '  +0.02%  '
Set-TextCell 'E28' '  -4.69%  '
Set-TextCell 'E29' '  -1.14%  '
Set-TextCell 'D30' '6.88'
Set-TextCell 'E30' '  -2.15%  '
Set-TextCell 'D31' '26.52'
Set-TextCell 'E32' '  +1.74%  '
Set-TextCell 'E33' '  +0.07%  '
Set-TextCell 'E34' '  +2.42%  '
Set-TextCell 'E35' '  -0.33%  '
Set-TextCell 'D36' '5.60'
Set-TextCell 'E36' '  -0.14%  '
Set-TextCell 'E37' '  -2.24%  '
Set-TextCell 'D38' '1.99'
Set-TextCell 'E38' '  +0.27%  '
Set-TextCell 'E39' '  -0.22%  '
Set-TextCell 'D40' '8.52'
Set-TextCell 'E40' '  -0.33%  '
Set-TextCell 'D41' '41.98'
Set-TextCell 'E41' '  +4.33%  '
Set-TextCell 'E42' '  -3.10%  '
Set-TextCell 'D43' '0.0343'
Set-TextCell 'E43' '  -0.71%  '
Set-TextCell 'D44' '2.687.71'
Set-TextCell 'E44' '  -0.39%  '
Set-TextCell 'D45' '133.60'
Set-TextCell 'E45' '  +2.05%  '
Set-TextCell 'D46' '361.36'
Set-TextCell 'E46' '  -3.40%  '
Set-TextCell 'E47' '  +0.07%  '
Set-TextCell 'D48' '23.47'
Set-TextCell 'E48' '  -1.70%  '
Set-TextCell 'E49' '  -1.16%  '
Set-TextCell 'E50' '  -1.10%  '
Set-TextCell 'E51' '  -0.98%  '
